$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.127881588408715
$ws.Range("C2").Value = 0.04240448674262143
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 0.8245905782990299

$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 15.68806981981553

$ws.Range("B4").Value = 0.6753301551942219
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 26.21740644021617
$ws.Range("E4").Value = 645.3272768299601
$ws.Range("G4").Value = 673.8878080086386

$ws.Range("B5").Value = 3.230985683306322
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 689428.5527653177
$ws.Range("E5").Value = 645.3272768299601
$ws.Range("G5").Value = 690078.7788224142

$ws.Range("B6").Value = 1.459612070389937
$ws.Range("C6").Value = 1.667794583268128
$ws.Range("D6").Value = 0.1575252929769615
$ws.Range("E6").Value = 0.496779210170732
$ws.Range("G6").Value = 3.781711156805759

$ws.Range("B7").Value = 0.01514828764759746
$ws.Range("C7").Value = 0.3127903958511391
$ws.Range("D7").Value = 3.900430680208489
$ws.Range("E7").Value = 8.660232485948974
$ws.Range("G7").Value = 12.8886018496562

$ws.Range("B8").Value = 0.04763786555579896
$ws.Range("C8").Value = 0.00007097389502863649
$ws.Range("D8").Value = 3.900430680208489
$ws.Range("E8").Value = 0.496779210170732
$ws.Range("G8").Value = 4.444918729830048

$ws.Range("B9").Value = 0.01514828764759746
$ws.Range("C9").Value = 0.04240448674262143
$ws.Range("D9").Value = 3.900430680208489
$ws.Range("E9").Value = 645.3272768299601
$ws.Range("G9").Value = 649.2852602845587

$ws.Range("B10").Value = 3.230985683306322
$ws.Range("C10").Value = 1.667794583268128
$ws.Range("D10").Value = 0.8054896365839992
$ws.Range("E10").Value = 0.496779210170732
$ws.Range("G10").Value = 6.201049113329182

$ws.Range("B11").Value = 0.127881588408715
$ws.Range("C11").Value = 0.04240448674262143
$ws.Range("D11").Value = 3.900430680208489
$ws.Range("E11").Value = 0.496779210170732
$ws.Range("G11").Value = 4.567495965530558
